$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 16
$srcRow = 15

# Copy formatting from the row above so the new row's styles (date format
# on column A, etc.) match the rest of the table.
$ws.Range("A" + $srcRow + ":N" + $srcRow).Copy()
$ws.Range("A" + $row + ":N" + $row).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42622.890416666669
$ws.Cells.Item($row, 2).Value = 18
$ws.Cells.Item($row, 3).Value = 62
$ws.Cells.Item($row, 4).Value = 33
$ws.Cells.Item($row, 5).Value = 62
$ws.Cells.Item($row, 6).Value = 22
$ws.Cells.Item($row, 7).Value = 25113
$ws.Cells.Item($row, 8).Value = 20148
$ws.Cells.Item($row, 9).Value = 3216
$ws.Cells.Item($row, 10).Value = 415
$ws.Cells.Item($row, 11).Value = 223
$ws.Cells.Item($row, 12).Value = 51
$ws.Cells.Item($row, 13).Value = 15
$ws.Cells.Item($row, 14).Value = "Noun"
